$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (e.g. "60.130.10") that must stay text, not be
# reinterpreted as a number -- format the whole data range as Text first,
# same as the source workbook already treats it, before writing new values.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '60.178.81'
$ws.Range('E2').Value = '  +2.43%  '
$ws.Range('D3').Value = '2.551.07'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '540.22'
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('D6').Value = '144.25'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Value = '0.572'
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '2.569.78'
$ws.Range('E9').Value = '  +2.56%  '
$ws.Range('E10').Value = '  +1.95%  '
$ws.Range('D11').Value = '0.161'
$ws.Range('E11').Value = '  +1.88%  '
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('D13').Value = '0.362'
$ws.Range('E13').Value = '  +3.49%  '
$ws.Range('D14').Value = '2.998.25'
$ws.Range('E14').Value = '  +1.58%  '
$ws.Range('D15').Value = '24.07'
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('D16').Value = '60.130.19'
$ws.Range('E17').Value = '  +3.91%  '
$ws.Range('D18').Value = '2.545.55'
$ws.Range('E18').Value = '  +1.47%  '
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('E20').Value = '  +1.59%  '
$ws.Range('D21').Value = '327.04'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '5.98'
$ws.Range('E22').Value = '  +4.38%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = '63.58'
$ws.Range('E24').Value = '  +4.60%  '
$ws.Range('E26').Value = '  +4.09%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  +4.37%  '
$ws.Range('D29').Value = '7.06'
$ws.Range('E29').Value = '  +3.08%  '
$ws.Range('D30').Value = '0.0₃0795'
$ws.Range('E30').Value = '  +4.30%  '
$ws.Range('E31').Value = '  +2.30%  '
$ws.Range('E32').Value = '  -3.77%  '
$ws.Range('D33').Value = '165.82'
$ws.Range('E33').Value = '  +5.58%  '
$ws.Range('E34').Value = '  +5.30%  '
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').Value = '18.75'
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('E38').Value = '  +2.26%  '
$ws.Range('E39').Value = '  +0.86%  '
$ws.Range('D40').Value = '5.62'
$ws.Range('E40').Value = '  -5.12%  '
$ws.Range('D41').Value = '301.14'
$ws.Range('E41').Value = '  -2.26%  '
$ws.Range('D42').Value = '3.73'
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('E43').Value = '  +5.91%  '
$ws.Range('E44').Value = '  +3.20%  '
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('D47').Value = '127.43'
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '18.96'
$ws.Range('E48').Value = '  +2.25%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.0939'
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('E51').Value = '  +1.34%  '
